# "Novas alterações com alura"
# - Lowercase the existing tipo_pista values (Dupla/Múltipla/Simples -> dupla/múltipla/simples)
# - Add a new row with the value "n/a"
# - Leave selection on the newly added cell

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Normalize existing values to lowercase (keep the trailing padding spaces intact)
$ws.Range("A2").Value = "dupla               "
$ws.Range("A3").Value = "múltipla            "
$ws.Range("A4").Value = "simples             "

# Append the new "n/a" entry on the next available row
$ws.Range("A5").Value = "n/a"

# Match the selection move onto the newly added cell
$ws.Range("A5").Select()
